$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2577.6
$ws.Range("I18").Value = 2577.6
$ws.Range("K18").Value = 2577.6
$ws.Range("M18").Value = -2293.6

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H121").Value = 3727.4
$ws.Range("J121").Value = 3727.4
$ws.Range("L121").Value = 11182.2
$ws.Range("N121").Value = -14676.2

$ws.Range("H127").Value = 1002652.2
$ws.Range("I127").Value = 1377710.6
$ws.Range("J127").Value = 2496.3333
$ws.Range("K127").Value = 4133131.8
$ws.Range("L127").Value = 7488.999899999999
$ws.Range("M127").Value = -4128171.8
$ws.Range("N127").Value = -17408.9999

$ws.Range("H137").Value = 1017.02704
$ws.Range("I137").Value = 947.21875
$ws.Range("J137").Value = 1463.8
$ws.Range("K137").Value = 2841.65625
$ws.Range("L137").Value = 4391.4
$ws.Range("M137").Value = -291.65625
$ws.Range("N137").Value = -9491.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 803.6142599999999
$ws.Range("I32").Value = 642.082
$ws.Range("K32").Value = 642.082
$ws.Range("M32").Value = -355.082

$ws.Range("I61").Value = 3356.5557
$ws.Range("K61").Value = 3356.5557
$ws.Range("M61").Value = -3144.5557

$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

$ws.Range("H122").Value = 2922.3125
$ws.Range("I122").Value = 2831.9092
$ws.Range("J122").Value = 3121.2
$ws.Range("K122").Value = 8495.7276
$ws.Range("L122").Value = 9363.599999999999
$ws.Range("M122").Value = -6045.7276
$ws.Range("N122").Value = -14263.6

$ws.Range("I136").Value = 3356.5557
$ws.Range("K136").Value = 10069.6671
$ws.Range("M136").Value = -7519.667099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5082.3477
$ws.Range("I86").Value = 5005.647
$ws.Range("K86").Value = 5005.647
$ws.Range("M86").Value = -3882.647

$ws.Range("H89").Value = 5082.3477
$ws.Range("I89").Value = 5005.647
$ws.Range("K89").Value = 25028.235
$ws.Range("M89").Value = -19412.235

$ws.Range("H105").Value = 3550.7896
$ws.Range("I105").Value = 2551
$ws.Range("J105").Value = 6350.2
$ws.Range("K105").Value = 2551
$ws.Range("L105").Value = 6350.2
$ws.Range("M105").Value = -804
$ws.Range("N105").Value = -9844.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2730.4167
$ws.Range("I16").Value = 2347
$ws.Range("J16").Value = 3497.25
$ws.Range("K16").Value = 2347
$ws.Range("L16").Value = 3497.25
$ws.Range("M16").Value = -2060
$ws.Range("N16").Value = -4071.25

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H31").Value = 307347.47
$ws.Range("J31").Value = 3546.4783
$ws.Range("L31").Value = 3546.4783
$ws.Range("N31").Value = -4136.478300000001

$ws.Range("H34").Value = 307347.47
$ws.Range("J34").Value = 3546.4783
$ws.Range("L34").Value = 3546.4783
$ws.Range("N34").Value = -3950.4783

$ws.Range("H62").Value = 4658.846
$ws.Range("J62").Value = 4756.857
$ws.Range("L62").Value = 4756.857
$ws.Range("N62").Value = -6004.857

$ws.Range("H65").Value = 4658.846
$ws.Range("J65").Value = 4756.857
$ws.Range("L65").Value = 23784.285
$ws.Range("N65").Value = -30024.285

$ws.Range("H113").Value = 2730.4167
$ws.Range("I113").Value = 2347
$ws.Range("J113").Value = 3497.25
$ws.Range("K113").Value = 2347
$ws.Range("L113").Value = 3497.25
$ws.Range("M113").Value = -177
$ws.Range("N113").Value = -7837.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 98.94118
$ws.Range("I2").Value = 36.125
$ws.Range("J2").Value = 154.77777
$ws.Range("K2").Value = 216.75
$ws.Range("L2").Value = 928.66662
$ws.Range("M2").Value = -103.75
$ws.Range("N2").Value = -1154.66662

$ws.Range("H11").Value = 71428710
$ws.Range("I11").Value = 142.375
$ws.Range("K11").Value = 427.125
$ws.Range("M11").Value = -287.125

$ws.Range("H12").Value = 341.63635
$ws.Range("J12").Value = 417.33334
$ws.Range("L12").Value = 1252.00002
$ws.Range("N12").Value = -1598.00002

$ws.Range("H131").Value = 1340.3529
$ws.Range("J131").Value = 1375
$ws.Range("L131").Value = 4125
$ws.Range("N131").Value = -14205

$ws.Range("H140").Value = 7833
$ws.Range("I140").Value = 6499
$ws.Range("K140").Value = 19497
$ws.Range("M140").Value = -14317

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4809576.5
$ws.Range("I2").Value = 5611151.5
$ws.Range("K2").Value = 5611151.5
$ws.Range("M2").Value = -5611038.5

$ws.Range("H80").Value = 2947.8667
$ws.Range("I80").Value = 2196
$ws.Range("K80").Value = 2196
$ws.Range("M80").Value = -1198

$ws.Range("H83").Value = 2947.8667
$ws.Range("I83").Value = 2196
$ws.Range("K83").Value = 10980
$ws.Range("M83").Value = -5988

$ws.Range("H100").Value = 46043.75
$ws.Range("J100").Value = 46043.75
$ws.Range("L100").Value = 46043.75
$ws.Range("N100").Value = -48207.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4215.5835
$ws.Range("J7").Value = 4441.143
$ws.Range("L7").Value = 4441.143
$ws.Range("N7").Value = -4665.143

$ws.Range("H22").Value = 2426.75
$ws.Range("I22").Value = 1888.375
$ws.Range("J22").Value = 3234.3125
$ws.Range("K22").Value = 1888.375
$ws.Range("L22").Value = 3234.3125
$ws.Range("M22").Value = -1593.375
$ws.Range("N22").Value = -3824.3125

$ws.Range("H27").Value = 2426.75
$ws.Range("I27").Value = 1888.375
$ws.Range("J27").Value = 3234.3125
$ws.Range("K27").Value = 1888.375
$ws.Range("L27").Value = 3234.3125
$ws.Range("M27").Value = -1781.375
$ws.Range("N27").Value = -3448.3125

$ws.Range("H30").Value = 1250
$ws.Range("I30").Value = 1250
$ws.Range("K30").Value = 1250
$ws.Range("M30").Value = -1142

$ws.Range("H40").Value = 4193.6924
$ws.Range("I40").Value = 3960.2727
$ws.Range("K40").Value = 3960.2727
$ws.Range("M40").Value = -3824.2727

$ws.Range("H61").Value = 11188.923
$ws.Range("I61").Value = 11769
$ws.Range("K61").Value = 11769
$ws.Range("M61").Value = -11567

$ws.Range("H82").Value = 3883.1667
$ws.Range("I82").Value = 3859.8
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 3859.8
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -3498.8
$ws.Range("N82").Value = -4722

$ws.Range("H85").Value = 3883.1667
$ws.Range("I85").Value = 3859.8
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 3859.8
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -2611.8
$ws.Range("N85").Value = -6496

$ws.Range("H93").Value = 2653.889
$ws.Range("J93").Value = 6264.5
$ws.Range("L93").Value = 6264.5
$ws.Range("N93").Value = -8760.5

$ws.Range("H113").Value = 11188.923
$ws.Range("I113").Value = 11769
$ws.Range("K113").Value = 11769
$ws.Range("M113").Value = -9599

$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws.Range("H126").Value = 4215.5835
$ws.Range("J126").Value = 4441.143
$ws.Range("L126").Value = 13323.429
$ws.Range("N126").Value = -18263.429

$ws.Range("H130").Value = 45000
$ws.Range("J130").Value = 45000
$ws.Range("L130").Value = 45000
$ws.Range("N130").Value = -55040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 503.6154
$ws.Range("I107").Value = 462.25
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1386.75
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 533.25
$ws.Range("N107").Value = -6840
